# TrialsSetup.xlsx update — 2026-01-23 12:00
# The "Days remaining" column (B) for three trials ticked down by one day
# since the last refresh:
#   Row 9  (ALLEGRETTO-LTE (B7981028)): 1  -> 0
#   Row 11 (REJOICE (MK-5909-003)):     28 -> 27
#   Row 14 (REMASTER (CLOU)):           48 -> 47

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 0
$ws.Range("B11").Value = 27
$ws.Range("B14").Value = 47
